$d = $word.ActiveDocument

# 1) Merge the split run "Súper Administrador tiene los mi" + "smos datos que Usuario."
#    into a single run reading "Súper Administrador tiene los mismos datos que Usuario."
#    (the old _GoBack bookmark that used to sit between the two runs is swallowed by the
#    replace and gets recreated from scratch below, in its new home).
$find = $d.Content.Find
$find.Execute("Súper Administrador tiene los mi" + "smos datos que Usuario.", `
              $true, $false, $false, $false, $false, $true, 1, $false, `
              "Súper Administrador tiene los mismos datos que Usuario.", 2) | Out-Null

# 2) Locate the paragraph holding that sentence so a new list paragraph can be
#    inserted right after it.
$targetIdx = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*mismos datos que Usuario*") {
        $targetIdx = $idx
    }
}

$target = $d.Paragraphs.Item($targetIdx)
$target.Range.InsertParagraphAfter()

# 3) Fill the new paragraph (it already inherited the "Prrafodelista" / numPr list
#    formatting from the paragraph it split off from). A trailing sentinel character
#    is appended along with the real sentence so the _GoBack bookmark below can be
#    anchored at an *interior* text position instead of sitting exactly on the
#    paragraph-mark boundary (the COM host mis-resolves degenerate ranges created
#    exactly at a paragraph mark). The sentinel is stripped off immediately after.
$newPara = $d.Paragraphs.Item($targetIdx + 1)
$newRange = $newPara.Range
$newRange.Collapse(1)
$newRange.InsertAfter("Falta ver tema de los retorno de los metodos~")

$newPara = $d.Paragraphs.Item($targetIdx + 1)
$sentinelPos = $newPara.Range.End - 2
$bmRange = $d.Range($sentinelPos, $sentinelPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$sentinelRange = $d.Range($sentinelPos, $sentinelPos + 1)
$sentinelRange.Delete()
